# Actualización automática 2025-11-28 17:30:09
# Updates November sales figures and recalculated dependent totals/percentages.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": PORCELANATO sales for row 6 (client) ---
$wsVentasPorGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasPorGrupo.Range("M6").Value = 24190.45

# --- Sheet "VENTA MENSUAL": noviembre column (F) ---
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F6").Value = 24929.52
$wsVentaMensual.Range("F26").Value = 58861.39

# --- Sheet "CUMPLIMIENTO MENSUAL": PORCELANATO row (12) and TOTAL row (14) ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D12").Value = 58122.32
$wsCumplimiento.Range("E12").Value = -23421.32
$wsCumplimiento.Range("F12").Value = 1.674946543327282

$wsCumplimiento.Range("D14").Value = 58861.39
$wsCumplimiento.Range("E14").Value = -18083.64941051808
$wsCumplimiento.Range("F14").Value = 1.443468646106953
